# Auto-generated edit script: updates F/G column numeric values
# (scrape counters: 'want-to-go' counts and min ticket prices) per commit diff.
$wb = $excel.ActiveWorkbook

# 展览 (Worksheets.Item(1))
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 3461
$ws.Range("F4").Value = 366
$ws.Range("F5").Value = 8107
$ws.Range("F8").Value = 2093
$ws.Range("F9").Value = 63
$ws.Range("F11").Value = 536
$ws.Range("F13").Value = 1
$ws.Range("F14").Value = 1058
$ws.Range("F18").Value = 1136
$ws.Range("F19").Value = 1
$ws.Range("F20").Value = 716
$ws.Range("F24").Value = 4615
$ws.Range("F25").Value = 100
$ws.Range("F26").Value = 48913
$ws.Range("F27").Value = 3977
$ws.Range("F29").Value = 984
$ws.Range("F30").Value = 735
$ws.Range("F31").Value = 53
$ws.Range("F32").Value = 66
$ws.Range("F35").Value = 562
$ws.Range("F36").Value = 188
$ws.Range("F38").Value = 561
$ws.Range("F39").Value = 824
$ws.Range("F40").Value = 981
$ws.Range("F41").Value = 117
$ws.Range("F42").Value = 152
$ws.Range("F43").Value = 1046
$ws.Range("F45").Value = 88
$ws.Range("F46").Value = 1
$ws.Range("F47").Value = 80
$ws.Range("F48").Value = 21
$ws.Range("F49").Value = 2449
# 演出 (Worksheets.Item(2))
$ws = $wb.Worksheets.Item(2)
$ws.Range("G12").Value = 180
$ws.Range("F16").Value = 80
$ws.Range("G16").Value = 180
$ws.Range("F20").Value = 7282
$ws.Range("F28").Value = 102
# 本地生活 (Worksheets.Item(3))
$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 2197
$ws.Range("F5").Value = 1472
$ws.Range("F7").Value = 625
$ws.Range("F8").Value = 2313
$ws.Range("F9").Value = 9236
$ws.Range("F10").Value = 1499
$ws.Range("F11").Value = 145
$ws.Range("F12").Value = 56
# 全部类型 (Worksheets.Item(4))
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 3461
$ws.Range("F3").Value = 2197
$ws.Range("F4").Value = 8107
$ws.Range("F5").Value = 1472
$ws.Range("F6").Value = 625
$ws.Range("F7").Value = 1499
$ws.Range("F8").Value = 145
$ws.Range("F10").Value = 63
$ws.Range("F11").Value = 536
$ws.Range("F13").Value = 1058
$ws.Range("F19").Value = 1136
$ws.Range("F22").Value = 4615
$ws.Range("F23").Value = 100
$ws.Range("G24").Value = 180
$ws.Range("F26").Value = 3977
$ws.Range("F28").Value = 984
$ws.Range("F29").Value = 735
$ws.Range("F30").Value = 54
$ws.Range("F31").Value = 66
$ws.Range("F33").Value = 562
$ws.Range("F35").Value = 188
$ws.Range("F36").Value = 561
$ws.Range("F37").Value = 824
$ws.Range("F39").Value = 117
$ws.Range("F40").Value = 152
$ws.Range("F41").Value = 1046
$ws.Range("F44").Value = 88
$ws.Range("F46").Value = 80
$ws.Range("F47").Value = 21
